# Scheduled runner update: refresh market-price-derived Leve profit columns
# (currentAveragePrice / NQ / HQ, LevePrice NQ/HQ, LeveProfit NQ/HQ) across sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 2591.1428  # ALC!H132
$ws.Cells.Item(132, 9).Value = 2329.6924  # ALC!I132
$ws.Cells.Item(132, 11).Value = 6989.0772  # ALC!K132
$ws.Cells.Item(132, 13).Value = -4459.0772  # ALC!M132

$ws.Cells.Item(137, 8).Value = 2615.3845  # ALC!H137
$ws.Cells.Item(137, 9).Value = 2135.3333  # ALC!I137
$ws.Cells.Item(137, 11).Value = 6405.999899999999  # ALC!K137
$ws.Cells.Item(137, 13).Value = -3855.999899999999  # ALC!M137

$ws.Cells.Item(138, 8).Value = 3301.5454  # ALC!H138
$ws.Cells.Item(138, 9).Value = 1148.5  # ALC!I138
$ws.Cells.Item(138, 10).Value = 3780  # ALC!J138
$ws.Cells.Item(138, 11).Value = 3445.5  # ALC!K138
$ws.Cells.Item(138, 12).Value = 11340  # ALC!L138
$ws.Cells.Item(138, 13).Value = 1694.5  # ALC!M138
$ws.Cells.Item(138, 14).Value = -21620  # ALC!N138

$ws.Cells.Item(141, 8).Value = 5526.8125  # ALC!H141
$ws.Cells.Item(141, 9).Value = 5228.8  # ALC!I141
$ws.Cells.Item(141, 11).Value = 15686.4  # ALC!K141
$ws.Cells.Item(141, 13).Value = -10506.4  # ALC!M141

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 277.5  # ARM!H2
$ws.Cells.Item(2, 9).Value = 263  # ARM!I2
$ws.Cells.Item(2, 11).Value = 263  # ARM!K2
$ws.Cells.Item(2, 13).Value = -150  # ARM!M2

$ws.Cells.Item(45, 8).Value = 1264.4375  # ARM!H45
$ws.Cells.Item(45, 9).Value = 1209.5714  # ARM!I45
$ws.Cells.Item(45, 11).Value = 1209.5714  # ARM!K45
$ws.Cells.Item(45, 13).Value = -832.5714  # ARM!M45

$ws.Cells.Item(76, 8).Value = 0  # ARM!H76
$ws.Cells.Item(76, 10).Value = 0  # ARM!J76
$ws.Cells.Item(76, 12).Value = 0  # ARM!L76
$ws.Cells.Item(76, 14).ClearContents()  # ARM!N76

$ws.Cells.Item(79, 8).Value = 0  # ARM!H79
$ws.Cells.Item(79, 10).Value = 0  # ARM!J79
$ws.Cells.Item(79, 12).Value = 0  # ARM!L79
$ws.Cells.Item(79, 14).ClearContents()  # ARM!N79

$ws.Cells.Item(80, 8).Value = 36244.25  # ARM!H80
$ws.Cells.Item(80, 9).Value = 10000  # ARM!I80
$ws.Cells.Item(80, 10).Value = 39993.43  # ARM!J80
$ws.Cells.Item(80, 11).Value = 10000  # ARM!K80
$ws.Cells.Item(80, 12).Value = 39993.43  # ARM!L80
$ws.Cells.Item(80, 13).Value = -9002  # ARM!M80
$ws.Cells.Item(80, 14).Value = -41989.43  # ARM!N80

$ws.Cells.Item(83, 8).Value = 36244.25  # ARM!H83
$ws.Cells.Item(83, 9).Value = 10000  # ARM!I83
$ws.Cells.Item(83, 10).Value = 39993.43  # ARM!J83
$ws.Cells.Item(83, 11).Value = 30000  # ARM!K83
$ws.Cells.Item(83, 12).Value = 119980.29  # ARM!L83
$ws.Cells.Item(83, 13).Value = -25008  # ARM!M83
$ws.Cells.Item(83, 14).Value = -129964.29  # ARM!N83

$ws.Cells.Item(116, 8).Value = 277.5  # ARM!H116
$ws.Cells.Item(116, 9).Value = 263  # ARM!I116
$ws.Cells.Item(116, 11).Value = 263  # ARM!K116
$ws.Cells.Item(116, 13).Value = 2031  # ARM!M116

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 277.5  # BSM!H3
$ws.Cells.Item(3, 9).Value = 263  # BSM!I3
$ws.Cells.Item(3, 11).Value = 263  # BSM!K3
$ws.Cells.Item(3, 13).Value = -149  # BSM!M3

$ws.Cells.Item(134, 8).Value = 3488.2903  # BSM!H134
$ws.Cells.Item(134, 10).Value = 3100  # BSM!J134
$ws.Cells.Item(134, 12).Value = 9300  # BSM!L134
$ws.Cells.Item(134, 14).Value = -14370  # BSM!N134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 3075.6924  # CRP!H58
$ws.Cells.Item(58, 9).Value = 1999.6666  # CRP!I58
$ws.Cells.Item(58, 11).Value = 1999.6666  # CRP!K58
$ws.Cells.Item(58, 13).Value = -1796.6666  # CRP!M58

$ws.Cells.Item(74, 8).Value = 38268.535  # CRP!H74
$ws.Cells.Item(74, 10).Value = 38268.535  # CRP!J74
$ws.Cells.Item(74, 12).Value = 38268.535  # CRP!L74
$ws.Cells.Item(74, 14).Value = -40016.535  # CRP!N74

$ws.Cells.Item(77, 8).Value = 38268.535  # CRP!H77
$ws.Cells.Item(77, 10).Value = 38268.535  # CRP!J77
$ws.Cells.Item(77, 12).Value = 114805.605  # CRP!L77
$ws.Cells.Item(77, 14).Value = -123541.605  # CRP!N77

$ws.Cells.Item(132, 8).Value = 1837.9375  # CRP!H132
$ws.Cells.Item(132, 10).Value = 2193  # CRP!J132
$ws.Cells.Item(132, 12).Value = 6579  # CRP!L132
$ws.Cells.Item(132, 14).Value = -11639  # CRP!N132

$ws.Cells.Item(134, 8).Value = 3995.6667  # CRP!H134
$ws.Cells.Item(134, 9).Value = 3995.5  # CRP!I134
$ws.Cells.Item(134, 11).Value = 11986.5  # CRP!K134
$ws.Cells.Item(134, 13).Value = -9451.5  # CRP!M134

$ws.Cells.Item(136, 8).Value = 3075.6924  # CRP!H136
$ws.Cells.Item(136, 9).Value = 1999.6666  # CRP!I136
$ws.Cells.Item(136, 11).Value = 5998.9998  # CRP!K136
$ws.Cells.Item(136, 13).Value = -3448.9998  # CRP!M136

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 2011.421  # CUL!H131
$ws.Cells.Item(131, 10).Value = 2699.4  # CUL!J131
$ws.Cells.Item(131, 12).Value = 8098.200000000001  # CUL!L131
$ws.Cells.Item(131, 14).Value = -18178.2  # CUL!N131

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 4978.8335  # GSM!H80
$ws.Cells.Item(80, 9).Value = 3495  # GSM!I80
$ws.Cells.Item(80, 10).Value = 5275.6  # GSM!J80
$ws.Cells.Item(80, 11).Value = 3495  # GSM!K80
$ws.Cells.Item(80, 12).Value = 5275.6  # GSM!L80
$ws.Cells.Item(80, 13).Value = -2497  # GSM!M80
$ws.Cells.Item(80, 14).Value = -7271.6  # GSM!N80

$ws.Cells.Item(83, 8).Value = 4978.8335  # GSM!H83
$ws.Cells.Item(83, 9).Value = 3495  # GSM!I83
$ws.Cells.Item(83, 10).Value = 5275.6  # GSM!J83
$ws.Cells.Item(83, 11).Value = 17475  # GSM!K83
$ws.Cells.Item(83, 12).Value = 26378  # GSM!L83
$ws.Cells.Item(83, 13).Value = -12483  # GSM!M83
$ws.Cells.Item(83, 14).Value = -36362  # GSM!N83

$ws.Cells.Item(102, 8).Value = 1368  # GSM!H102
$ws.Cells.Item(102, 9).Value = 1186.6666  # GSM!I102
$ws.Cells.Item(102, 10).Value = 3000  # GSM!J102
$ws.Cells.Item(102, 11).Value = 1186.6666  # GSM!K102
$ws.Cells.Item(102, 12).Value = 3000  # GSM!L102
$ws.Cells.Item(102, 13).Value = 435.3334  # GSM!M102
$ws.Cells.Item(102, 14).Value = -6244  # GSM!N102

$ws.Cells.Item(132, 8).Value = 811.4706  # GSM!H132
$ws.Cells.Item(132, 9).Value = 811.4706  # GSM!I132
$ws.Cells.Item(132, 11).Value = 2434.4118  # GSM!K132
$ws.Cells.Item(132, 13).Value = 95.58820000000014  # GSM!M132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(81, 8).Value = 0  # LTW!H81
$ws.Cells.Item(81, 10).Value = 0  # LTW!J81
$ws.Cells.Item(81, 12).Value = 0  # LTW!L81
$ws.Cells.Item(81, 14).ClearContents()  # LTW!N81

$ws.Cells.Item(84, 8).Value = 0  # LTW!H84
$ws.Cells.Item(84, 10).Value = 0  # LTW!J84
$ws.Cells.Item(84, 12).Value = 0  # LTW!L84
$ws.Cells.Item(84, 14).ClearContents()  # LTW!N84

$ws.Cells.Item(100, 8).Value = 2137.6  # LTW!H100
$ws.Cells.Item(100, 9).Value = 2219.5557  # LTW!I100
$ws.Cells.Item(100, 11).Value = 2219.5557  # LTW!K100
$ws.Cells.Item(100, 13).Value = -1678.5557  # LTW!M100

$ws.Cells.Item(136, 8).Value = 2872.5715  # LTW!H136
$ws.Cells.Item(136, 10).Value = 3432.3333  # LTW!J136
$ws.Cells.Item(136, 12).Value = 10296.9999  # LTW!L136
$ws.Cells.Item(136, 14).Value = -15396.9999  # LTW!N136

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(32, 8).Value = 46000  # WVR!H32
$ws.Cells.Item(32, 9).Value = 0  # WVR!I32
$ws.Cells.Item(32, 10).Value = 46000  # WVR!J32
$ws.Cells.Item(32, 11).Value = 0  # WVR!K32
$ws.Cells.Item(32, 12).Value = 46000  # WVR!L32
$ws.Cells.Item(32, 13).ClearContents()  # WVR!M32
$ws.Cells.Item(32, 14).Value = -46634  # WVR!N32

$ws.Cells.Item(68, 8).Value = 70000  # WVR!H68
$ws.Cells.Item(68, 10).Value = 70000  # WVR!J68
$ws.Cells.Item(68, 12).Value = 70000  # WVR!L68
$ws.Cells.Item(68, 14).Value = -71622  # WVR!N68

$ws.Cells.Item(71, 8).Value = 70000  # WVR!H71
$ws.Cells.Item(71, 10).Value = 70000  # WVR!J71
$ws.Cells.Item(71, 12).Value = 210000  # WVR!L71
$ws.Cells.Item(71, 14).Value = -218112  # WVR!N71

$ws.Cells.Item(80, 8).Value = 0  # WVR!H80
$ws.Cells.Item(80, 10).Value = 0  # WVR!J80
$ws.Cells.Item(80, 12).Value = 0  # WVR!L80
$ws.Cells.Item(80, 14).ClearContents()  # WVR!N80

$ws.Cells.Item(82, 8).Value = 57499  # WVR!H82
$ws.Cells.Item(82, 10).Value = 57499  # WVR!J82
$ws.Cells.Item(82, 12).Value = 57499  # WVR!L82
$ws.Cells.Item(82, 14).Value = -58265  # WVR!N82

$ws.Cells.Item(83, 8).Value = 0  # WVR!H83
$ws.Cells.Item(83, 10).Value = 0  # WVR!J83
$ws.Cells.Item(83, 12).Value = 0  # WVR!L83
$ws.Cells.Item(83, 14).ClearContents()  # WVR!N83

$ws.Cells.Item(85, 8).Value = 57499  # WVR!H85
$ws.Cells.Item(85, 10).Value = 57499  # WVR!J85
$ws.Cells.Item(85, 12).Value = 57499  # WVR!L85
$ws.Cells.Item(85, 14).Value = -60151  # WVR!N85

$ws.Cells.Item(126, 8).Value = 2758.7144  # WVR!H126
$ws.Cells.Item(126, 9).Value = 2077.75  # WVR!I126
$ws.Cells.Item(126, 10).Value = 3666.6667  # WVR!J126
$ws.Cells.Item(126, 11).Value = 6233.25  # WVR!K126
$ws.Cells.Item(126, 12).Value = 11000.0001  # WVR!L126
$ws.Cells.Item(126, 13).Value = -3763.25  # WVR!M126
$ws.Cells.Item(126, 14).Value = -15940.0001  # WVR!N126

$ws.Cells.Item(132, 8).Value = 747.6  # WVR!H132
$ws.Cells.Item(132, 9).Value = 539.4286  # WVR!I132
$ws.Cells.Item(132, 11).Value = 1618.2858  # WVR!K132
$ws.Cells.Item(132, 13).Value = 911.7142000000001  # WVR!M132

$ws.Cells.Item(136, 8).Value = 2231.9473  # WVR!H136
$ws.Cells.Item(136, 9).Value = 2189.2778  # WVR!I136
$ws.Cells.Item(136, 11).Value = 6567.8334  # WVR!K136
$ws.Cells.Item(136, 13).Value = -4017.8334  # WVR!M136

